$d = $word.ActiveDocument

# 1) "...des effets bonus/malus !" -> "...des effets bonus."
$d.Content.Find.Execute("bonus/malus !", $true, $false, $false, $false, $false, $true, 1, $false, "bonus.", 2) | Out-Null

# 2) "...sur les effets bonus/malus. Les enfants..." -> "...sur les effets bonus. Les enfants..."
$d.Content.Find.Execute("bonus/malus. Les enfants", $true, $false, $false, $false, $false, $true, 1, $false, "bonus. Les enfants", 2) | Out-Null

# 3) "...leur nombre de bonus/malus." -> "...leur nombre de bonus."
$d.Content.Find.Execute("leur nombre de bonus/malus.", $true, $false, $false, $false, $false, $true, 1, $false, "leur nombre de bonus.", 2) | Out-Null

# 4) Delete the whole "TODO : 1) tester..." paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("TODO")) {
        $p.Range.Delete()
        break
    }
}

# 5) Add a sentence at the end of the "hexagones verts et gris" paragraph.
$d.Content.Find.Execute("hexagones verts et gris comme sur la figure 1.", $true, $false, $false, $false, $false, $true, 1, $false, "hexagones verts et gris comme sur la figure 1. Les hexagones gris représentent les zones d’en-but. TODO : colorer la zone d’en-but à la couleur des équipiers qui la défende.", 2) | Out-Null

# 6) Rework the "jeton blanc" sentence.
$d.Content.Find.Execute("Le jeton blanc représente le disque volant que lance et attrapent les pièces des joueurs représentées par des empilements de 1, 2 ou", $true, $false, $false, $false, $false, $true, 1, $false, "Le disque volant représenté par le jeton blanc que lance et attrapent les équipiers bleus ou rouges  qui sont représentés par des empilements de 1, 2 ou", 2) | Out-Null

# 7) "Chaque pièce peut recevoir" -> "Chaque équipier peut recevoir"
$d.Content.Find.Execute("Chaque pièce peut recevoir", $true, $false, $false, $false, $false, $true, 1, $false, "Chaque équipier peut recevoir", 2) | Out-Null

# 8) "2 compteurs affichent" -> "2 dés affichent"
$d.Content.Find.Execute("Enfin, 2 compteurs affichent les scores de chaque joueur.", $true, $false, $false, $false, $false, $true, 1, $false, "Enfin, 2 dés affichent les scores de chaque joueur.", 2) | Out-Null

# 9) "Aucun bonus n'est posé sur les pièces." -> "Aucun bonus n'est posé sur les équipiers."
$d.Content.Find.Execute("Aucun bonus n’est posé sur les pièces. Chaque joueur pioche", $true, $false, $false, $false, $false, $true, 1, $false, "Aucun bonus n’est posé sur les équipiers. Chaque joueur pioche", 2) | Out-Null

# 10) Insert a new list item "Non ! Car des stratégies de barrage..." before
#     "Simplifier les règles d'engagement,".
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Simplifier les règles d’engagement")) {
        $target = $p
        break
    }
}
$target.Range.InsertParagraphBefore()
$target.Range.Text = "Non ! Car des stratégies de barrage en occupant à deux les case seraient possibles !"
$target.Range.ListFormat.ListLevelNumber = 2

# 11) Merge "Précision : " and "Pour marquer un point..." (no visible text change).
$d.Content.Find.Execute("Pour marquer un point, l’équipier dans la zone d'en but doit recevoir le disque et en obtenir la possession.", $true, $false, $false, $false, $false, $true, 1, $false, "Pour marquer un point, l’équipier dans la zone d'en but doit recevoir le disque et en obtenir la possession.", 2) | Out-Null

# 12) Update the SAVEDATE field text in the header.
$sec = $d.Sections.Item(1)
$header = $sec.Headers.Item(1)
$header.Range.Find.Execute("2022-0809-2055", $true, $false, $false, $false, $false, $true, 1, $false, "2022-0810-1000", 2) | Out-Null
